# "can check for votes"
# Add a new data row (row 6) to the form-responses sheet containing a
# vote/election code (123456) in column B, matching the style already
# used by the other rows in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vote entry.
$ws.Range("B6").Value = 123456

# Give B6 the same number format / font as the existing B2:B5 entries
# (copy formats only, so the numeric value just set above is preserved).
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)

# Make sure the row renders at the sheet's normal row height.
$ws.Rows.Item(6).RowHeight = 15.75

# Leave the selection where the author left it when they saved the file.
$null = $ws.Range("D12").Select()
